$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Set the Title value (B5) to match the Name value (B4): "DESCNonQualifiant"
$ws.Range("B5").Value = "DESCNonQualifiant"

# Update the Date value (B8) to reflect the new generation timestamp
$ws.Range("B8").Value = "2025-07-17T14:35:50+00:00"
